$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 3.3
$ws.Range("J3").Value = 2.6
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("W3").Value = 6
$ws.Range("Y3").Value = 9
$ws.Range("AA3").Value = 17
$ws.Range("AC3").Value = 7.5
$ws.Range("AG3").Value = 10
$ws.Range("AO3").Value = 10
$ws.Range("AP3").Value = 23
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("BA3").Value = 126
# Row 5
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 4.2
$ws.Range("K5").Value = 1.91
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("S5").Value = 1.62
$ws.Range("T5").Value = 2.2
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 8
$ws.Range("Y5").Value = 10
$ws.Range("AC5").Value = 6
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 81
$ws.Range("AG5").Value = 8.5
$ws.Range("AP5").Value = 29
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 301
$ws.Range("AT5").Value = 2.2
$ws.Range("AU5").Value = 9.5
$ws.Range("AV5").Value = 81
$ws.Range("AX5").Value = 26
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 101
$ws.Range("BA5").Value = 151
# Row 7
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 4.65
$ws.Range("I7").Value = 7
$ws.Range("L7").Value = 6.2
$ws.Range("Q7").Value = 1.55
$ws.Range("R7").Value = 2.15
$ws.Range("S7").Value = 1.29
$ws.Range("T7").Value = 3.42
$ws.Range("V7").Value = 1.82
$ws.Range("AD7").Value = 9.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 21
$ws.Range("AI7").Value = 22
$ws.Range("AJ7").Value = 150
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 600
$ws.Range("AO7").Value = 6.1
$ws.Range("AT7").Value = 3.25
$ws.Range("AU7").Value = 8
$ws.Range("AV7").Value = 70
$ws.Range("AX7").Value = 37
# Row 8
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 1.85
$ws.Range("K8").Value = 2.15
$ws.Range("P8").Value = 3.2
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.9
$ws.Range("T8").Value = 2.99
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 55
$ws.Range("AC8").Value = 10.75
$ws.Range("AD8").Value = 6.8
$ws.Range("AG8").Value = 7.7
$ws.Range("AK8").Value = 14
$ws.Range("AU8").Value = 7.1
# Row 9
$ws.Range("O9").Value = 1.17
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 1.6
$ws.Range("R9").Value = 2.3
# Row 15
$ws.Range("G15").Value = 2.3
$ws.Range("I15").Value = 3.3
$ws.Range("J15").Value = 3.1
$ws.Range("L15").Value = 4
$ws.Range("X15").Value = 10
$ws.Range("Z15").Value = 21
$ws.Range("AA15").Value = 21
$ws.Range("AC15").Value = 7.5
$ws.Range("AG15").Value = 8.5
$ws.Range("AH15").Value = 15
$ws.Range("AI15").Value = 13
$ws.Range("AJ15").Value = 34
$ws.Range("AK15").Value = 29
$ws.Range("AM15").Value = 351
$ws.Range("AO15").Value = 13
$ws.Range("AP15").Value = 26
$ws.Range("AQ15").Value = 41
$ws.Range("AR15").Value = 67
$ws.Range("AS15").Value = 201
$ws.Range("AW15").Value = 5
$ws.Range("AX15").Value = 19
$ws.Range("AZ15").Value = 67
$ws.Range("BA15").Value = 101
# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("I16").Value = 3.25
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("Q16").Value = 2.08
$ws.Range("R16").Value = 1.73
$ws.Range("U16").Value = 1.83
$ws.Range("V16").Value = 1.83
$ws.Range("W16").Value = 7.5
$ws.Range("X16").Value = 10
$ws.Range("AS16").Value = 201
$ws.Range("AY16").Value = 29
# Row 19
$ws.Range("G19").Value = 2.2
$ws.Range("I19").Value = 3.3
$ws.Range("J19").Value = 2.88
$ws.Range("L19").Value = 3.75
$ws.Range("W19").Value = 8
$ws.Range("X19").Value = 11
$ws.Range("Z19").Value = 21
$ws.Range("AA19").Value = 19
$ws.Range("AE19").Value = 13
$ws.Range("AI19").Value = 12
$ws.Range("AJ19").Value = 34
$ws.Range("AK19").Value = 26
$ws.Range("AM19").Value = 201
$ws.Range("AN19").Value = 4.33
$ws.Range("AO19").Value = 12
$ws.Range("AW19").Value = 5
$ws.Range("AY19").Value = 26
$ws.Range("AZ19").Value = 51
# Row 20
$ws.Range("G20").Value = 2.15
$ws.Range("I20").Value = 3.8
$ws.Range("J20").Value = 2.88
$ws.Range("X20").Value = 9.5
$ws.Range("Z20").Value = 19
$ws.Range("AG20").Value = 10
$ws.Range("AK20").Value = 34
$ws.Range("AO20").Value = 12
$ws.Range("AP20").Value = 23
$ws.Range("AW20").Value = 5.5
$ws.Range("AX20").Value = 21
# Row 22
$ws.Range("Q22").Value = 2.1
$ws.Range("R22").Value = 1.7
# Row 25
$ws.Range("G25").Value = 2.1
$ws.Range("H25").Value = 3.7
$ws.Range("I25").Value = 2.85
$ws.Range("K25").Value = 2.32
$ws.Range("L25").Value = 3.35
$ws.Range("M25").Value = 1.03
$ws.Range("N25").Value = 9
$ws.Range("O25").Value = 1.19
$ws.Range("P25").Value = 4.2
$ws.Range("Q25").Value = 1.57
$ws.Range("R25").Value = 2.25
$ws.Range("S25").Value = 1.3
$ws.Range("T25").Value = 3.2
$ws.Range("U25").Value = 1.53
$ws.Range("V25").Value = 2.35
$ws.Range("W25").Value = 10.75
$ws.Range("X25").Value = 12.5
$ws.Range("Y25").Value = 9
$ws.Range("Z25").Value = 21
$ws.Range("AA25").Value = 15
$ws.Range("AB25").Value = 21
$ws.Range("AC25").Value = 9
$ws.Range("AD25").Value = 7.7
$ws.Range("AE25").Value = 12
$ws.Range("AF25").Value = 40
$ws.Range("AG25").Value = 12.5
$ws.Range("AI25").Value = 10.5
$ws.Range("AJ25").Value = 35
$ws.Range("AK25").Value = 21
$ws.Range("AL25").Value = 25
$ws.Range("AM25").Value = 250
$ws.Range("AN25").Value = 4.3
$ws.Range("AP25").Value = 16.5
$ws.Range("AQ25").Value = 37
$ws.Range("AR25").Value = 60
$ws.Range("AS25").Value = 175
$ws.Range("AT25").Value = 3.2
$ws.Range("AU25").Value = 6.7
$ws.Range("AV25").Value = 50
$ws.Range("AW25").Value = 5.1
$ws.Range("AX25").Value = 15
$ws.Range("AY25").Value = 19.5
$ws.Range("AZ25").Value = 60
$ws.Range("BA25").Value = 80
$ws.Range("BB25").Value = 200
# Row 27
$ws.Range("G27").Value = 2.18
$ws.Range("H27").Value = 2.95
$ws.Range("I27").Value = 3.45
$ws.Range("J27").Value = 2.82
$ws.Range("K27").Value = 1.93
$ws.Range("L27").Value = 4
$ws.Range("M27").Value = 1.04
$ws.Range("N27").Value = 6.35
$ws.Range("O27").Value = 1.44
$ws.Range("P27").Value = 2.42
$ws.Range("Q27").Value = 2.25
$ws.Range("R27").Value = 1.5
$ws.Range("S27").Value = 1.47
$ws.Range("T27").Value = 2.32
$ws.Range("U27").Value = 1.93
$ws.Range("V27").Value = 1.7
$ws.Range("W27").Value = 6
$ws.Range("X27").Value = 9.5
$ws.Range("Z27").Value = 21
$ws.Range("AB27").Value = 37
$ws.Range("AC27").Value = 6.9
$ws.Range("AD27").Value = 5.8
$ws.Range("AE27").Value = 16.5
$ws.Range("AF27").Value = 100
$ws.Range("AG27").Value = 8.25
$ws.Range("AH27").Value = 17.5
$ws.Range("AI27").Value = 12
$ws.Range("AJ27").Value = 50
$ws.Range("AK27").Value = 35
$ws.Range("AL27").Value = 50
$ws.Range("AM27").Value = 1000
$ws.Range("AO27").Value = 11.75
$ws.Range("AP27").Value = 22
$ws.Range("AT27").Value = 2.3
$ws.Range("AU27").Value = 7.3
$ws.Range("AV27").Value = 75
$ws.Range("AW27").Value = 5.1
$ws.Range("AX27").Value = 19.5
$ws.Range("AY27").Value = 28
$ws.Range("AZ27").Value = 110
$ws.Range("BA27").Value = 150
$ws.Range("BB27").Value = 400
